# Applies the "completed regression testing" update to the Testing Matrix workbook.
# - Marks a large batch of test cells as completed (checkmark) on the "GUI" sheet
# - Clears the stale "Two tests: complet and cancel" note in J4
# - Adds a new "Dispose" test row (row 60) with its own completed checkmark
# - Leaves final selection on F56:G56 / view scrolled to A42, matching the saved file

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GUI")

$checkA = [char]0x2714
$checkB = [char]0xFE0F
$check = [string]$checkA + [string]$checkB   # "✔️" heavy check mark + variation selector-16

# Rows 4-7 : "Groups Testing" sub-table -> all of B:G fully checked off
"4","5","6","7" | ForEach-Object {
    $ws.Range("B" + $_ + ":G" + $_).Value = $check
}

# Remove the stale note that used to sit in J4
$ws.Range("J4").ClearContents()

# Rows 13-16 : "By Size" sub-table -> all of B:G fully checked off
"13","14","15","16" | ForEach-Object {
    $ws.Range("B" + $_ + ":G" + $_).Value = $check
}
# G15 previously had no cell/formatting at all (row only went to F) -- match the
# green centered style used by the rest of the "checked" cells in that column
$ws.Range("G15").Font.Color = $ws.Range("F15").Font.Color
$ws.Range("G15").HorizontalAlignment = $ws.Range("F15").HorizontalAlignment

# Rows 20-23 : "By Exposure" sub-table -> all of B:G fully checked off
"20","21","22","23" | ForEach-Object {
    $ws.Range("B" + $_ + ":G" + $_).Value = $check
}
# G22 likewise needs the shared "checked" cell style applied
$ws.Range("G22").Font.Color = $ws.Range("F22").Font.Color
$ws.Range("G22").HorizontalAlignment = $ws.Range("F22").HorizontalAlignment

# Rows 27-30 : "By Temperature" sub-table -> all of B:G fully checked off
"27","28","29","30" | ForEach-Object {
    $ws.Range("B" + $_ + ":G" + $_).Value = $check
}
# G29 likewise needs the shared "checked" cell style applied
$ws.Range("G29").Font.Color = $ws.Range("F29").Font.Color
$ws.Range("G29").HorizontalAlignment = $ws.Range("F29").HorizontalAlignment

# Row 34 : "Size and Exposure" / "No Precal" -> everything but E34
$ws.Range("B34").Value = $check
$ws.Range("C34").Value = $check
$ws.Range("D34").Value = $check
$ws.Range("F34").Value = $check
$ws.Range("G34").Value = $check

# Row 41 : "Size and Temperature" / "No Precal" -> only B41
$ws.Range("B41").Value = $check

# Row 42 : "Size and Temperature" / "Pedestal" -> D42 and G42
$ws.Range("D42").Value = $check
$ws.Range("G42").Value = $check

# Row 44 : "Size and Temperature" / "Fixed file" -> F44 and G44
$ws.Range("F44").Value = $check
$ws.Range("G44").Value = $check

# Row 48 : "Exposure and Temperature" / "No Precal" -> only B48
$ws.Range("B48").Value = $check

# Row 56 : "All 3" / "No Precal" -> B56, F56, G56
$ws.Range("B56").Value = $check
$ws.Range("F56").Value = $check
$ws.Range("G56").Value = $check

# Row 57 : "All 3" / "Pedestal" -> only C57
$ws.Range("C57").Value = $check

# New row 60 : "Dispose" test, with its own completed checkmark in C60
$ws.Range("A60").Value = "Dispose"
$ws.Range("C60").Value = $check
# C60 is a brand new cell -- apply the same green centered "checked" style used elsewhere
$ws.Range("C60").Font.Color = $ws.Range("B56").Font.Color
$ws.Range("C60").HorizontalAlignment = $ws.Range("B56").HorizontalAlignment

# Match the saved view state: scrolled so row 42 is at the top, with F56:G56 selected
$ws.Activate()
$ws.Range("A42").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 42
$ws.Range("F56:G56").Select() | Out-Null
